$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Il18 -> Il18r1 -> Resolving-Mac (D2 changes from ECs to Resolving-Mac)
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 0.2184343333333333
$ws.Range("H2").Value = 0.655303
$ws.Range("I2").Value = 0.008416673064019609
$ws.Range("J2").Value = 0.00841667306401961
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2710126666666666
$ws.Range("N2").Value = 0.8130379999999999
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.05919847116822221
$ws.Range("R2").Value = 0.5327862405139999
$ws.Range("S2").Value = 0.008416673064019609
$ws.Range("T2").Value = 0.00841667306401961

# Row 3: A3 changes from ECs to FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.199828
$ws.Range("H3").Value = 12.599484
$ws.Range("I3").Value = 0.1618270290283213
$ws.Range("J3").Value = 0.1618270290283213
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2710126666666666
$ws.Range("N3").Value = 0.8130379999999999
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.138206585821333
$ws.Range("R3").Value = 10.243859272392
$ws.Range("S3").Value = 0.1618270290283213
$ws.Range("T3").Value = 0.1618270290283213

# Row 4: A4 changes from FAPs to MuSCs, D4 changes from ECs to Resolving-Mac
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 3.307112333333333
$ws.Range("H4").Value = 9.921336999999999
$ws.Range("I4").Value = 0.1274290669918512
$ws.Range("J4").Value = 0.1274290669918513
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2710126666666666
$ws.Range("N4").Value = 0.8130379999999999
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.8962693324228886
$ws.Range("R4").Value = 8.066423991805999
$ws.Range("S4").Value = 0.1274290669918512
$ws.Range("T4").Value = 0.1274290669918513

# Row 5: A5 changes from FAPs to Resolving-Mac
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 18.22719966666667
$ws.Range("H5").Value = 54.681599
$ws.Range("I5").Value = 0.7023272309158078
$ws.Range("J5").Value = 0.7023272309158078
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2710126666666666
$ws.Range("N5").Value = 0.8130379999999999
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 4.93980198752911
$ws.Range("R5").Value = 44.458217887762
$ws.Range("S5").Value = 0.7023272309158078
$ws.Range("T5").Value = 0.7023272309158078

# Rows 6-9 are removed entirely (the data previously there is gone now).
$ws.Rows("6:9").Delete()
